$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.118.16"
$ws.Range("E2").Value = "  -0.92%  "
$ws.Range("D3").Value = "2.574.34"
$ws.Range("E3").Value = "  -2.69%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "588.20"
$ws.Range("E5").Value = "  -3.32%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "149.60"
$ws.Range("E6").Value = "  +1.06%  "
$ws.Range("E7").Value = "  +0.19%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.584"
$ws.Range("E8").Value = "  -0.59%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.109"
$ws.Range("E9").Value = "  +0.47%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.64"
$ws.Range("E10").Value = "  +1.59%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.382"
$ws.Range("E11").Value = "  +0.32%  "
$ws.Range("E12").Value = "  -0.67%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "27.33"
$ws.Range("E13").Value = "  -0.42%  "
$ws.Range("D14").Value = "3.045.31"
$ws.Range("E14").Value = "  -2.01%  "
$ws.Range("D15").Value = "63.004.74"
$ws.Range("E15").Value = "  -0.79%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000155"
$ws.Range("E16").Value = "  +5.42%  "
$ws.Range("D17").Value = "2.578.65"
$ws.Range("E17").Value = "  -2.67%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.28"
$ws.Range("E18").Value = "  +4.88%  "
$ws.Range("E19").Value = "  +3.03%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "344.73"
$ws.Range("E20").Value = "  -0.25%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.86"
$ws.Range("E21").Value = "  -0.55%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "67.41"
$ws.Range("E23").Value = "  +1.61%  "
$ws.Range("E24").Value = "  +3.58%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.20"
$ws.Range("E25").Value = "  +1.17%  "
$ws.Range("E26").Value = "  -1.55%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "548.58"
$ws.Range("E27").Value = "  -2.65%  "
$ws.Range("E28").Value = "  +0.84%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.96"
$ws.Range("E29").Value = "  -1.39%  "
$ws.Range("E30").Value = "  -1.31%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.03"
$ws.Range("E31").Value = "  -0.31%  "
$ws.Range("D32").Value = "0.0₃0842"
$ws.Range("E32").Value = "  -1.10%  "
$ws.Range("E33").Value = "  -1.37%  "
$ws.Range("E34").Value = "  -2.46%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "167.27"
$ws.Range("E35").Value = "  -0.91%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.410"
$ws.Range("E36").Value = "  +1.46%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.998"
$ws.Range("E37").Value = "  -0.16%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.46"
$ws.Range("E38").Value = "  +1.75%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.91"
$ws.Range("E39").Value = "  -0.51%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  +0.11%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "166.21"
$ws.Range("E41").Value = "  +0.73%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "39.55"
$ws.Range("E42").Value = "  -1.13%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.92"
$ws.Range("E43").Value = "  +3.44%  "
$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "22.36"
$ws.Range("E44").Value = "  +1.45%  "
$ws.Range("B45").Value = "Hedera"
$ws.Range("C45").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0580"
$ws.Range("E45").Value = "  +2.08%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.628"
$ws.Range("E46").Value = "  -0.34%  "
$ws.Range("E47").Value = "  +2.42%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.03"
$ws.Range("E48").Value = "  +1.18%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0960"
$ws.Range("E49").Value = "  +0.23%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "18.98"
$ws.Range("E50").Value = "  +0.56%  "
$ws.Range("D51").Value = "0.0₆0229"
$ws.Range("E51").Value = "  +16.37%  "
